$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '29.355.42'
$ws.Range('E2').Value = '  +0.06%  '

# Row 3
Set-TextValue 'D3' '1.839.18'
$ws.Range('E3').Value = '  -0.25%  '

# Row 4
Set-TextValue 'D4' '0.9995'
$ws.Range('E4').Value = '  +0.14%  '

# Row 5
Set-TextValue 'D5' '238.72'
$ws.Range('E5').Value = '  -0.49%  '

# Row 6
$ws.Range('E6').Value = '  -0.14%  '

# Row 7
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
Set-TextValue 'D8' '0.07403'
$ws.Range('E8').Value = '  -1.02%  '

# Row 9
Set-TextValue 'D9' '0.2887'
$ws.Range('E9').Value = '  -0.36%  '

# Row 10
Set-TextValue 'D10' '24.88'
$ws.Range('E10').Value = '  +1.70%  '

# Row 11
Set-TextValue 'D11' '0.07714'
$ws.Range('E11').Value = '  -0.29%  '

# Row 12
Set-TextValue 'D12' '1.837.32'
$ws.Range('E12').Value = '  -0.32%  '

# Row 13
Set-TextValue 'D13' '4.964'
$ws.Range('E13').Value = '  -0.35%  '

# Row 14
Set-TextValue 'D14' '0.6731'
$ws.Range('E14').Value = '  -0.95%  '

# Row 15
Set-TextValue 'D15' '0.00001023'
$ws.Range('E15').Value = '  -3.31%  '

# Row 16
Set-TextValue 'D16' '81.65'
$ws.Range('E16').Value = '  -0.34%  '

# Row 17
Set-TextValue 'D17' '6.199'
$ws.Range('E17').Value = '  +0.32%  '

# Row 18
Set-TextValue 'D18' '29.447.75'
$ws.Range('E18').Value = '  +0.35%  '

# Row 19
Set-TextValue 'D19' '232.66'
$ws.Range('E19').Value = '  +1.49%  '

# Row 20
Set-TextValue 'D20' '12.29'
$ws.Range('E20').Value = '  -0.18%  '

# Row 21
$ws.Range('E21').Value = '  +0.16%  '

# Row 22
Set-TextValue 'D22' '7.283'
$ws.Range('E22').Value = '  -2.95%  '

# Row 23
Set-TextValue 'D23' '1.001'
$ws.Range('E23').Value = '  +0.25%  '

# Row 24
Set-TextValue 'D24' '157.89'
$ws.Range('E24').Value = '  -0.26%  '

# Row 25
Set-TextValue 'D25' '8.478'
$ws.Range('E25').Value = '  +0.64%  '

# Row 26
Set-TextValue 'D26' '0.1343'
$ws.Range('E26').Value = '  -1.78%  '

# Row 27
Set-TextValue 'D27' '17.28'
$ws.Range('E27').Value = '  -1.36%  '

# Row 28
Set-TextValue 'D28' '0.07251'
$ws.Range('E28').Value = '  +10.10%  '

# Row 29
Set-TextValue 'D29' '1.476'

# Row 30
Set-TextValue 'D30' '1.477'
$ws.Range('E30').Value = '  -0.39%  '

# Row 31
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D31' '4.022'
$ws.Range('E31').Value = '  -1.64%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D32' '4.027'
$ws.Range('E32').Value = '  -1.88%  '

# Row 33
Set-TextValue 'D33' '1.812'
$ws.Range('E33').Value = '  -0.74%  '

# Row 34
Set-TextValue 'D34' '1.138'
$ws.Range('E34').Value = '  -0.11%  '

# Row 35
Set-TextValue 'D35' '0.6961'
$ws.Range('E35').Value = '  +0.24%  '

# Row 36
Set-TextValue 'D36' '2.572'
$ws.Range('E36').Value = '  -0.27%  '

# Row 37
Set-TextValue 'D37' '0.01836'
$ws.Range('E37').Value = '  +0.03%  '

# Row 38
Set-TextValue 'D38' '6.898'
$ws.Range('E38').Value = '  +1.53%  '

# Row 39
Set-TextValue 'D39' '2.814'
$ws.Range('E39').Value = '  -0.66%  '

# Row 40
Set-TextValue 'D40' '1.232.87'
$ws.Range('E40').Value = '  -2.48%  '

# Row 41
Set-TextValue 'D41' '0.9457'
$ws.Range('E41').Value = '  +3.13%  '

# Row 42
$ws.Range('E42').Value = '  +0.22%  '

# Row 43
Set-TextValue 'D43' '2.015.92'
$ws.Range('E43').Value = '  +0.67%  '

# Row 44
Set-TextValue 'D44' '100.73'
$ws.Range('E44').Value = '  -0.49%  '

# Row 45
Set-TextValue 'D45' '65.26'
$ws.Range('E45').Value = '  -1.27%  '

# Row 46
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D46' '0.00000000117'
$ws.Range('E46').Value = '  -1.02%  '

# Row 47
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D47' '1.706'
$ws.Range('E47').Value = '  -1.57%  '

# Row 48
Set-TextValue 'D48' '6.930'
$ws.Range('E48').Value = '  -1.93%  '

# Row 49
Set-TextValue 'D49' '8.841'
$ws.Range('E49').Value = '  -1.35%  '

# Row 50
Set-TextValue 'D50' '0.3890'
$ws.Range('E50').Value = '  -1.45%  '

# Row 51
Set-TextValue 'D51' '0.1128'
$ws.Range('E51').Value = '  -2.86%  '
